$p = $ppt.ActivePresentation
$s = $p.Slides.Item(16)
$tbl = $s.Shapes.Item(3).Table
$tbl.ApplyStyle("{CC981B3C-168A-43DC-AF3F-2BCA7A1978A2}")
